$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('20160405_164857', 2209.889, 'trim "space" and ",", convert to lower, remove multiple spaces, convert unicode to ascii', '8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_1, first_character_ascii, first_character_digit', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 4000', 0.993333333333333, 0.99009900990099, '0 filters: ', 0.322916666666667),
    @('20160405_172547', 2221.802, 'trim "space" and ",", convert to lower, remove multiple spaces, convert unicode to ascii', '8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_1, first_character_ascii, first_character_digit', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 4000', 0.992, 0.993399339933993, '0 filters: ', 0.206185567010309),
    @('20160405_180248', 2294.271, 'trim "space" and ",", convert to lower, remove multiple spaces, convert unicode to ascii', '8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_1, first_character_ascii, first_character_digit', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 4000', 0.990666666666667, 0.99009900990099, '0 filters: ', 0.291666666666667),
    @('20160405_184103', 2319.464, 'trim "space" and ",", convert to lower, remove multiple spaces, convert unicode to ascii', '8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_1, first_character_ascii, first_character_digit', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 4000', 0.994, 0.99009900990099, '0 filters: ', 0.28125),
    @('20160405_191942', 2357.317, 'trim "space" and ",", convert to lower, remove multiple spaces, convert unicode to ascii', '8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_1, first_character_ascii, first_character_digit', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 4000', 0.992, 0.99009900990099, '0 filters: ', 0.239583333333333),
    @('20160406_081430', 3338.312, 'remove multiple spaces, convert unicode to ascii, trim "space" and ",", convert to lower', '8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_1, first_character_ascii, first_character_digit', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 6000', 0.992666666666667, 0.99009900990099, '0 filters: ', 0.28125),
    @('20160406_091009', 3454.685, 'remove multiple spaces, convert unicode to ascii, trim "space" and ",", convert to lower', '8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_1, first_character_ascii, first_character_digit', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 6000', 0.99, 0.993399339933993, '0 filters: ', 0.164948453608247),
    @('20160406_100743', 5036.042, 'remove multiple spaces, convert unicode to ascii, trim "space" and ",", convert to lower', '8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_1, first_character_ascii, first_character_digit', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 6000', 0.994, 0.99009900990099, '0 filters: ', 0.302083333333333),
    @('20160406_113139', 5760.941, 'remove multiple spaces, convert unicode to ascii, trim "space" and ",", convert to lower', '8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_1, first_character_ascii, first_character_digit', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 6000', 0.992, 0.993399339933993, '0 filters: ', 0.206185567010309),
    @('20160406_130740', 6939.185, 'remove multiple spaces, convert unicode to ascii, trim "space" and ",", convert to lower', '8 features: #ascii/(#ascii+#digit+#punctuation), #digit/(#ascii+#digit+#punctuation), %kwName, %kwAddress, %kwPhone, #max_digit_skip_0_1, first_character_ascii, first_character_digit', 'Neuron Network', '2 layers: [100-Sigmoid, 3-Softmax], learning_rate: 0.01, learning_rule: adagrad, n_iterator: 6000', 0.994, 0.99009900990099, '0 filters: ', 0.322916666666667),
)

$startRow = 22
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $ws.Cells.Item($r, 9).Value = $row[8]
    $ws.Cells.Item($r, 10).Value = $row[9]
}
